$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ordering of (id, speaker_variant) pairs for rows 2..22.
# This reflects an export with no "is_prefered" marker and no
# Levenshtein-distance based sort (column D is cleared for all rows).
$rows = @(
    @{ R = 2;  B = "#polinos";                 C = "Polinos" },
    @{ R = 3;  B = "#ceres";                   C = "Ceres" },
    @{ R = 4;  B = "#jaager";                  C = "Jaager" },
    @{ R = 5;  B = "#armida";                  C = "Armida" },
    @{ R = 6;  B = "#hubaldus";                C = "Hubaldus" },
    @{ R = 7;  B = "#reinout,-armida,-dares";  C = "Reinout, Armida, Dares" },
    @{ R = 8;  B = "#dares";                   C = "Dares" },
    @{ R = 9;  B = "#karel";                   C = "Karel" },
    @{ R = 10; B = "#bacchus";                 C = "Bacchus" },
    @{ R = 11; B = "#politionelle";            C = "Politionelle" },
    @{ R = 12; B = "#rfinout";                 C = "Rfinout" },
    @{ R = 13; B = "#dans";                    C = "Dans" },
    @{ R = 14; B = "#pedestaalen";             C = "Pedestaalen" },
    @{ R = 15; B = "#thisbe";                  C = "Thisbe" },
    @{ R = 16; B = "#lukvrouw";                C = "Lukvrouw" },
    @{ R = 17; B = "#reinout";                 C = "Reinout" },
    @{ R = 18; B = "#aurora";                  C = "Aurora" },
    @{ R = 19; B = "#kupido";                  C = "Kupido" },
    @{ R = 20; B = "#hydraot";                 C = "Hydraot" },
    @{ R = 21; B = "#filida";                  C = "Filida" },
    @{ R = 22; B = "#aap";                     C = "Aap" }
)

foreach ($row in $rows) {
    $r = $row.R
    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    # Clear the is_prefered column - no preferred marker in this export.
    $ws.Range("D$r").Value = ""
}
